$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 94

# Text columns (A-D). The interop layer auto-detects dates/times/numbers in
# plain strings and silently converts them (e.g. "2024-01-26" -> date serial,
# "16:13:19" -> time serial). Temporarily force a text number format while
# assigning the values, then restore "General" so the new row's formatting
# matches the rest of the sheet (which never had an explicit text format).
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2024-01-26"
$ws.Range("B$row").Value = "16:13:19"
$ws.Range("C$row").Value = "Friday"
$ws.Range("D$row").Value = "03"

$textRange.NumberFormat = "General"

# Numeric columns (E-T)
$ws.Range("E$row").Value = 136373
$ws.Range("F$row").Value = 141654
$ws.Range("G$row").Value = 171391
$ws.Range("H$row").Value = 149284
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 122406
$ws.Range("K$row").Value = 223801
$ws.Range("L$row").Value = 256683
$ws.Range("M$row").Value = 185266
$ws.Range("N$row").Value = 110033
$ws.Range("O$row").Value = 41365
$ws.Range("P$row").Value = 30820
$ws.Range("Q$row").Value = 73569
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 42325
$ws.Range("T$row").Value = -1
